# Update loading_percent values for rows 2-25 (corresponding to the "case with 380 kV" run).
# Columns A, D, J, L are untouched (they stay 0 / index values); B,C,E,F,G,H,I,K,M,N,O are updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "B" = 7.24749128347; "C" = 5.608244049269726; "E" = 13.21107360574126; "F" = 16.86991607391245; "G" = 22.26228912007355; "H" = 12.78597287451445; "I" = 18.26165178571049; "K" = 8.279231250906387; "M" = 12.87468222577749; "N" = 17.51817861921914; "O" = 18.57301426752911 }
    3 = @{ "B" = 6.939430588932619; "C" = 5.508998543217425; "E" = 12.99699615621783; "F" = 15.89584955866815; "G" = 22.35556290807613; "H" = 12.82998753553376; "I" = 18.34597976387236; "K" = 8.053272603469672; "M" = 12.69609018909863; "N" = 17.56672653108214; "O" = 18.64956194940254 }
    4 = @{ "B" = 6.743681868787177; "C" = 5.446572816577784; "E" = 12.86841301823543; "F" = 15.26997757108489; "G" = 22.42022031313576; "H" = 12.85879324644586; "I" = 18.40093858598185; "K" = 7.909838938330773; "M" = 12.58806153920596; "N" = 17.59807926552753; "O" = 18.70022514835245 }
    5 = @{ "B" = 6.662365131775553; "C" = 5.420779110183501; "E" = 12.81681008934248; "F" = 15.008197319934; "G" = 22.44841704270807; "H" = 12.87097993927804; "I" = 18.42413521125074; "K" = 7.850260391073549; "M" = 12.54450242142478; "N" = 17.61124501211892; "O" = 18.72179091954314 }
    6 = @{ "B" = 6.648772631830398; "C" = 5.416475248956991; "E" = 12.80829172266638; "F" = 14.96433081551589; "G" = 22.45321044635615; "H" = 12.873030608435; "I" = 18.42803535556852; "K" = 7.840300888677909; "M" = 12.53729902224255; "N" = 17.61345471289918; "O" = 18.72542744588489 }
    7 = @{ "B" = 6.742591310543382; "C" = 5.446226363394262; "E" = 12.86771376147947; "F" = 15.26647399323133; "G" = 22.42059311274921; "H" = 12.85895578543785; "I" = 18.40124818200862; "K" = 7.909039937035041; "M" = 12.58747213824792; "N" = 17.59825524602837; "O" = 18.70051226757538 }
    8 = @{ "B" = 7.142703951855546; "C" = 5.574344664412831; "E" = 13.13671198251954; "F" = 16.5399640634477; "G" = 22.29291134604323; "H" = 12.80077982209571; "I" = 18.29006815983861; "K" = 8.20232500285179; "M" = 12.81279856925573; "N" = 17.53459804242424; "O" = 18.59864719960339 }
    9 = @{ "B" = 7.90860067638415; "C" = 5.812954777893348; "E" = 13.68351923569379; "F" = 19.00274580682531; "G" = 22.10153450559298; "H" = 12.70080639525396; "I" = 18.09726016033795; "K" = 8.738130264027127; "M" = 13.26517514135763; "N" = 17.42197376302952; "O" = 18.4279909248318 }
    10 = @{ "B" = 8.477461341941; "C" = 5.979530148447665; "E" = 14.09246255652562; "F" = 20.67494806633232; "G" = 21.99741671228051; "H" = 12.63593206537272; "I" = 17.97094137021003; "K" = 9.105318635746718; "M" = 13.60072833824517; "N" = 17.3466077812066; "O" = 18.32040454815495 }
    11 = @{ "B" = 8.722544866668425; "C" = 6.053213525416847; "E" = 14.27911508099577; "F" = 21.3917225636224; "G" = 21.95807465993182; "H" = 12.60827612780918; "I" = 17.91679927155738; "K" = 9.266166324173833; "M" = 13.75340179532223; "N" = 17.3139108760767; "O" = 18.27533661029522 }
    12 = @{ "B" = 8.813365936647703; "C" = 6.080799090263749; "E" = 14.34980823899814; "F" = 21.65686569030329; "G" = 21.9443378862707; "H" = 12.59807005749376; "I" = 17.89677426324286; "K" = 9.326152684674295; "M" = 13.81116477962764; "N" = 17.3017567134239; "O" = 18.25882865928999 }
    13 = @{ "B" = 8.793894589022052; "C" = 6.074872373444434; "E" = 14.33458388044953; "F" = 21.60004134736742; "G" = 21.94724458447488; "H" = 12.60025626208609; "I" = 17.90106577474864; "K" = 9.313275100809735; "M" = 13.79872762491259; "N" = 17.30436422715889; "O" = 18.26235909104831 }
    14 = @{ "B" = 8.730056673251269; "C" = 6.055489441319761; "E" = 14.28493124618699; "F" = 21.4136618050453; "G" = 21.95692121280515; "H" = 12.60743112573565; "I" = 17.9151422342158; "K" = 9.271120118438967; "M" = 13.7581553266439; "N" = 17.31290639214814; "O" = 18.27396729476651 }
    15 = @{ "B" = 8.690694927200649; "C" = 6.043575135272079; "E" = 14.25451686574132; "F" = 21.29868154950795; "G" = 21.96299986979796; "H" = 12.61186065416483; "I" = 17.92382664463476; "K" = 9.245177820785116; "M" = 13.73329534918416; "N" = 17.31816831063625; "O" = 18.28115039965814 }
    16 = @{ "B" = 8.46116739614078; "C" = 5.974671283570161; "E" = 14.08027020589759; "F" = 20.62722412089977; "G" = 22.00014989514079; "H" = 12.63777676888422; "I" = 17.97454649860263; "K" = 9.094679411994944; "M" = 13.59074658650312; "N" = 17.34877647339952; "O" = 18.3234279002039 }
    17 = @{ "B" = 8.316838392942667; "C" = 5.931853473197072; "E" = 13.9734788965509; "F" = 20.20408069597325; "G" = 22.02500042415461; "H" = 12.65415060141313; "I" = 18.00651196396572; "K" = 9.000743761857228; "M" = 13.50326508111281; "N" = 17.36795958181628; "O" = 18.35035681431059 }
    18 = @{ "B" = 8.232535957598733; "C" = 5.907029720053374; "E" = 13.91211841702663; "F" = 19.95656407809801; "G" = 22.04004840003017; "H" = 12.66374307812367; "I" = 18.02521024255945; "K" = 8.946133766485847; "M" = 13.45295394641887; "N" = 17.37914266612553; "O" = 18.36621014248997 }
    19 = @{ "B" = 8.203771958375414; "C" = 5.898591644669039; "E" = 13.89135604764653; "F" = 19.87204792380568; "G" = 22.04527274726109; "H" = 12.66702093445162; "I" = 18.03159485007617; "K" = 8.927545098097415; "M" = 13.43592218465285; "N" = 17.38295476303356; "O" = 18.37164037656235 }
    20 = @{ "B" = 8.332335927657947; "C" = 5.93643191355738; "E" = 13.98484105757089; "F" = 20.24955283636154; "G" = 22.02227688618203; "H" = 12.65238950174192; "I" = 18.00307683487532; "K" = 9.010803724171344; "M" = 13.51257742982393; "N" = 17.36590204387254; "O" = 18.34745244884006 }
    21 = @{ "B" = 8.748861454689923; "C" = 6.061191395525111; "E" = 14.29951571075252; "F" = 21.46857628470577; "G" = 21.95404737900481; "H" = 12.60531645968095; "I" = 17.91099468030643; "K" = 9.283527361609185; "M" = 13.7700742146317; "N" = 17.31039118495117; "O" = 18.27054252168592 }
    22 = @{ "B" = 9.009499467358802; "C" = 6.140876132172293; "E" = 14.50520909653044; "F" = 22.22866616901552; "G" = 21.91622703526744; "H" = 12.57610555566992; "I" = 17.85359646655099; "K" = 9.456374022162457; "M" = 13.93804069813758; "N" = 17.27543698251399; "O" = 18.22353222791889 }
    23 = @{ "B" = 8.871456845287776; "C" = 6.098521392453357; "E" = 14.39544836366556; "F" = 21.82633154458858; "G" = 21.93579046366133; "H" = 12.59155383745824; "I" = 17.88397637628738; "K" = 9.364626186065001; "M" = 13.8484411742253; "N" = 17.2939716985776; "O" = 18.24832427027451 }
    24 = @{ "B" = 8.325333620331417; "C" = 5.934362644934428; "E" = 13.97970410878455; "F" = 20.22900810905287; "G" = 22.0235058275996; "H" = 12.65318513749025; "I" = 18.00462885712301; "K" = 9.006257497539199; "M" = 13.50836736522393; "N" = 17.3668317755771; "O" = 18.34876435525065 }
    25 = @{ "B" = 7.687113556947622; "C" = 5.749866607602653; "E" = 13.53400036660234; "F" = 18.34778573295695; "G" = 22.14693509030964; "H" = 12.72634401460281; "I" = 18.14672407513895; "K" = 8.597663827301737; "M" = 13.14200354154307; "N" = 17.45114128411849; "O" = 18.47103698206923 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowData = $newValues[$rowNum]
    foreach ($colLetter in $rowData.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowData[$colLetter]
    }
}
